# Generate Report for Handback
# Applies the "handback" localization-status update:
#  - Status text "In Translation" -> "Handed back: in sync with en-US"
#    on the Overview sheet and on each per-locale sheet's Status column.
#  - Fill in "Latest Target File" (hyperlinked) and "Latest Handback File"
#    columns, plus a "Latest Handback DateTime" stamp, for both rows on
#    the zh-cn and de-de sheets.
#  - Widen the columns that now hold the longer strings.

$wb = $excel.ActiveWorkbook

$statusOld = "In Translation"
$statusNew = "Handed back: in sync with en-US"

$mdFile1 = "5a4b14dc-a148-4efb-aec8-0eaf41633d0d.md"
$mdFile2 = "95ec2a14-f20d-4d17-8b45-4fb2989ecf7b.md"
$mdUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf24fb0ac128924d57fb50a5f08796ceb27acdbb/e2e/5a4b14dc-a148-4efb-aec8-0eaf41633d0d.md"
$mdUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf24fb0ac128924d57fb50a5f08796ceb27acdbb/e2e/95ec2a14-f20d-4d17-8b45-4fb2989ecf7b.md"

$zhHandback1 = "5a4b14dc-a148-4efb-aec8-0eaf41633d0d.f05a79eed66a0b810ce201f74c78542f643a48b9.zh-cn.xlf"
$zhHandback2 = "95ec2a14-f20d-4d17-8b45-4fb2989ecf7b.0efcec51733a9006c4285fe1bd08cf6a0ed876f6.zh-cn.xlf"
$deHandback1 = "5a4b14dc-a148-4efb-aec8-0eaf41633d0d.f05a79eed66a0b810ce201f74c78542f643a48b9.de-de.xlf"
$deHandback2 = "95ec2a14-f20d-4d17-8b45-4fb2989ecf7b.0efcec51733a9006c4285fe1bd08cf6a0ed876f6.de-de.xlf"

$zhHandbackDate = "2016-08-31 20:31:43"
$deHandbackDate = "2016-08-31 20:31:49"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsZh.Range("J2").Value = $zhHandback1
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdFile2)
$wsZh.Range("J3").Value = $zhHandback2
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsDe.Range("J2").Value = $deHandback1
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdFile2)
$wsDe.Range("J3").Value = $deHandback2
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
